$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q3" right after "总计" (position 2), pushing
#    the existing quarter sheets (2021-Q4, 2021-Q2, 2021-Q1, 2020-Q4) down by
#    one slot each.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($null, $totalSheet)
$newSheet.Name = "2022-Q3"

# Match the page margins used throughout the rest of the workbook (0.75in /
# 0.75in / 1in / 1in / 0.5in / 0.5in -- PageSetup units are points).
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data row (fund code / name / size / position stay text, rank is numeric --
# same convention as the other quarter sheets already in the workbook)
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'001541"
$newSheet.Range("B2").Style = "Normal"
$newSheet.Range("C2").Value = "汇添富民营新动力股票"
$newSheet.Range("D2").Value = "'2.21"
$newSheet.Range("D2").Style = "Normal"
$newSheet.Range("E2").Value = "'91.50"
$newSheet.Range("E2").Style = "Normal"
$newSheet.Range("F2").Value = "'3.29"
$newSheet.Range("F2").Style = "Normal"
$newSheet.Range("G2").Value = "'0.0727"
$newSheet.Range("G2").Style = "Normal"
$newSheet.Range("H2").Value = 8

# Re-use the header/index style (bold, centered, bordered) already defined on
# the "总计" sheet instead of re-declaring fonts/borders by hand.
$totalSheet.Range("A2").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$newSheet.Range("A2").PasteSpecial(-4122)
$newSheet.Range("A2").Value = 0

# ---------------------------------------------------------------------------
# 2. Update the "总计" summary sheet: push the existing quarters' figures
#    down one row and insert the new 2022-Q3 figures at the top of the table.
# ---------------------------------------------------------------------------
$oldB2 = $totalSheet.Range("B2").Value()
$oldC2 = $totalSheet.Range("C2").Value()
$oldD2 = $totalSheet.Range("D2").Value()
$oldB3 = $totalSheet.Range("B3").Value()
$oldC3 = $totalSheet.Range("C3").Value()
$oldD3 = $totalSheet.Range("D3").Value()
$oldB4 = $totalSheet.Range("B4").Value()
$oldC4 = $totalSheet.Range("C4").Value()
$oldD4 = $totalSheet.Range("D4").Value()
$oldB5 = $totalSheet.Range("B5").Value()
$oldC5 = $totalSheet.Range("C5").Value()
$oldD5 = $totalSheet.Range("D5").Value()

$totalSheet.Range("B6").Value = $oldB5
$totalSheet.Range("C6").Value = $oldC5
$totalSheet.Range("D6").Value = $oldD5

$totalSheet.Range("B5").Value = $oldB4
$totalSheet.Range("C5").Value = $oldC4
$totalSheet.Range("D5").Value = $oldD4

$totalSheet.Range("B4").Value = $oldB3
$totalSheet.Range("C4").Value = $oldC3
$totalSheet.Range("D4").Value = $oldD3

$totalSheet.Range("B3").Value = $oldB2
$totalSheet.Range("C3").Value = $oldC2
$totalSheet.Range("D3").Value = $oldD2

$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0.07000000000000001

# New row 6 (2020-Q4) needs the same row-index style as the rest of column A
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A5").Copy()
$totalSheet.Range("A6").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3. Restore the originally active sheet/tab (the last quarter sheet,
#    "2020-Q4") since adding a new sheet made it the active one.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Activate()
